$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Delete columns N:O (Recall Drowsy, Recall Non-Drowsy)
$ws.Range("N1:O1").EntireColumn.Delete() | Out-Null

# Step 2: Insert a new blank column at M to make room for Validation Accuracy
$ws.Range("M1").EntireColumn.Insert() | Out-Null

# Step 3: Update header labels
$ws.Range("L1").Value = "Train Accuracy"
$ws.Range("M1").Value = "Validation Accuracy"

# Step 4: Set Train Accuracy (L) and Validation Accuracy (M) values for rows 2-97
$lArr = New-Object 'object[,]' 96,1
$mArr = New-Object 'object[,]' 96,1
$lArr[0,0] = 0.9791666865348816
$mArr[0,0] = 1
$lArr[1,0] = 0.9708333611488342
$mArr[1,0] = 0.9833333492279053
$lArr[2,0] = 0.987500011920929
$mArr[2,0] = 0.8666666746139526
$lArr[3,0] = 0.9854166507720947
$mArr[3,0] = 1
$lArr[4,0] = 1
$mArr[4,0] = 0.9833333492279053
$lArr[5,0] = 0.9770833253860474
$mArr[5,0] = 0.9666666388511658
$lArr[6,0] = 0.949999988079071
$mArr[6,0] = 0.550000011920929
$lArr[7,0] = 0.96875
$mArr[7,0] = 0.9833333492279053
$lArr[8,0] = 0.9854166507720947
$mArr[8,0] = 0.9833333492279053
$lArr[9,0] = 0.9833333492279053
$mArr[9,0] = 1
$lArr[10,0] = 0.9729166626930237
$mArr[10,0] = 0.949999988079071
$lArr[11,0] = 0.9916666746139526
$mArr[11,0] = 1
$lArr[12,0] = 0.9729166626930237
$mArr[12,0] = 0.9666666388511658
$lArr[13,0] = 0.9916666746139526
$mArr[13,0] = 1
$lArr[14,0] = 0.9916666746139526
$mArr[14,0] = 0.8500000238418579
$lArr[15,0] = 0.9895833134651184
$mArr[15,0] = 0.9833333492279053
$lArr[16,0] = 0.9895833134651184
$mArr[16,0] = 0.9833333492279053
$lArr[17,0] = 0.987500011920929
$mArr[17,0] = 0.8666666746139526
$lArr[18,0] = 0.9666666388511658
$mArr[18,0] = 0.8999999761581421
$lArr[19,0] = 0.9750000238418579
$mArr[19,0] = 0.9666666388511658
$lArr[20,0] = 0.9791666865348816
$mArr[20,0] = 0.6166666746139526
$lArr[21,0] = 0.9750000238418579
$mArr[21,0] = 0.9833333492279053
$lArr[22,0] = 0.9895833134651184
$mArr[22,0] = 1
$lArr[23,0] = 0.9854166507720947
$mArr[23,0] = 0.949999988079071
$lArr[24,0] = 0.987500011920929
$mArr[24,0] = 0.949999988079071
$lArr[25,0] = 0.9770833253860474
$mArr[25,0] = 0.9833333492279053
$lArr[26,0] = 0.9937499761581421
$mArr[26,0] = 0.8833333253860474
$lArr[27,0] = 0.9895833134651184
$mArr[27,0] = 0.9666666388511658
$lArr[28,0] = 0.9854166507720947
$mArr[28,0] = 0.9833333492279053
$lArr[29,0] = 0.9916666746139526
$mArr[29,0] = 0.8833333253860474
$lArr[30,0] = 0.9750000238418579
$mArr[30,0] = 0.949999988079071
$lArr[31,0] = 0.9604166746139526
$mArr[31,0] = 0.800000011920929
$lArr[32,0] = 0.9833333492279053
$mArr[32,0] = 0.9333333373069763
$lArr[33,0] = 0.981249988079071
$mArr[33,0] = 0.949999988079071
$lArr[34,0] = 0.987500011920929
$mArr[34,0] = 0.9833333492279053
$lArr[35,0] = 0.9895833134651184
$mArr[35,0] = 0.8333333134651184
$lArr[36,0] = 0.9791666865348816
$mArr[36,0] = 0.8166666626930237
$lArr[37,0] = 0.9770833253860474
$mArr[37,0] = 0.8999999761581421
$lArr[38,0] = 0.9916666746139526
$mArr[38,0] = 0.8999999761581421
$lArr[39,0] = 0.9854166507720947
$mArr[39,0] = 0.9833333492279053
$lArr[40,0] = 0.9895833134651184
$mArr[40,0] = 0.9666666388511658
$lArr[41,0] = 0.9937499761581421
$mArr[41,0] = 0.8999999761581421
$lArr[42,0] = 0.9604166746139526
$mArr[42,0] = 0.800000011920929
$lArr[43,0] = 0.9770833253860474
$mArr[43,0] = 0.8666666746139526
$lArr[44,0] = 0.9854166507720947
$mArr[44,0] = 0.75
$lArr[45,0] = 0.9770833253860474
$mArr[45,0] = 0.8999999761581421
$lArr[46,0] = 0.9854166507720947
$mArr[46,0] = 0.9833333492279053
$lArr[47,0] = 0.9895833134651184
$mArr[47,0] = 0.9666666388511658
$lArr[48,0] = 0.9750000238418579
$mArr[48,0] = 0.9666666388511658
$lArr[49,0] = 0.9729166626930237
$mArr[49,0] = 1
$lArr[50,0] = 0.9916666746139526
$mArr[50,0] = 0.9833333492279053
$lArr[51,0] = 0.9833333492279053
$mArr[51,0] = 0.8999999761581421
$lArr[52,0] = 0.9791666865348816
$mArr[52,0] = 0.800000011920929
$lArr[53,0] = 0.9895833134651184
$mArr[53,0] = 0.9833333492279053
$lArr[54,0] = 0.9770833253860474
$mArr[54,0] = 0.949999988079071
$lArr[55,0] = 0.9833333492279053
$mArr[55,0] = 0.8999999761581421
$lArr[56,0] = 0.9916666746139526
$mArr[56,0] = 1
$lArr[57,0] = 0.9791666865348816
$mArr[57,0] = 0.9833333492279053
$lArr[58,0] = 0.9937499761581421
$mArr[58,0] = 0.9333333373069763
$lArr[59,0] = 0.9770833253860474
$mArr[59,0] = 1
$lArr[60,0] = 0.9708333611488342
$mArr[60,0] = 0.9166666865348816
$lArr[61,0] = 0.9854166507720947
$mArr[61,0] = 0.949999988079071
$lArr[62,0] = 0.9833333492279053
$mArr[62,0] = 0.6166666746139526
$lArr[63,0] = 0.9854166507720947
$mArr[63,0] = 0.9333333373069763
$lArr[64,0] = 0.9854166507720947
$mArr[64,0] = 0.9666666388511658
$lArr[65,0] = 0.9791666865348816
$mArr[65,0] = 0.8833333253860474
$lArr[66,0] = 0.9604166746139526
$mArr[66,0] = 0.9333333373069763
$lArr[67,0] = 0.9729166626930237
$mArr[67,0] = 0.9666666388511658
$lArr[68,0] = 0.9729166626930237
$mArr[68,0] = 0.9833333492279053
$lArr[69,0] = 0.9770833253860474
$mArr[69,0] = 0.8333333134651184
$lArr[70,0] = 0.987500011920929
$mArr[70,0] = 0.9833333492279053
$lArr[71,0] = 0.987500011920929
$mArr[71,0] = 0.8833333253860474
$lArr[72,0] = 0.981249988079071
$mArr[72,0] = 0.8166666626930237
$lArr[73,0] = 0.9916666746139526
$mArr[73,0] = 0.8666666746139526
$lArr[74,0] = 0.987500011920929
$mArr[74,0] = 0.8166666626930237
$lArr[75,0] = 0.9895833134651184
$mArr[75,0] = 0.8833333253860474
$lArr[76,0] = 0.9916666746139526
$mArr[76,0] = 0.9333333373069763
$lArr[77,0] = 0.987500011920929
$mArr[77,0] = 0.9833333492279053
$lArr[78,0] = 0.9645833373069763
$mArr[78,0] = 0.9666666388511658
$lArr[79,0] = 0.9645833373069763
$mArr[79,0] = 0.9166666865348816
$lArr[80,0] = 0.9791666865348816
$mArr[80,0] = 0.9833333492279053
$lArr[81,0] = 0.9770833253860474
$mArr[81,0] = 1
$lArr[82,0] = 0.9770833253860474
$mArr[82,0] = 0.8999999761581421
$lArr[83,0] = 0.9854166507720947
$mArr[83,0] = 0.9833333492279053
$lArr[84,0] = 0.9729166626930237
$mArr[84,0] = 0.9833333492279053
$lArr[85,0] = 0.9854166507720947
$mArr[85,0] = 0.9666666388511658
$lArr[86,0] = 0.9770833253860474
$mArr[86,0] = 0.6166666746139526
$lArr[87,0] = 0.981249988079071
$mArr[87,0] = 0.9166666865348816
$lArr[88,0] = 0.9770833253860474
$mArr[88,0] = 0.9333333373069763
$lArr[89,0] = 0.9895833134651184
$mArr[89,0] = 0.949999988079071
$lArr[90,0] = 0.9729166626930237
$mArr[90,0] = 0.9833333492279053
$lArr[91,0] = 0.9791666865348816
$mArr[91,0] = 0.9833333492279053
$lArr[92,0] = 0.9854166507720947
$mArr[92,0] = 0.949999988079071
$lArr[93,0] = 0.9770833253860474
$mArr[93,0] = 0.9833333492279053
$lArr[94,0] = 0.9895833134651184
$mArr[94,0] = 1
$lArr[95,0] = 0.9895833134651184
$mArr[95,0] = 0.9833333492279053

$ws.Range("L2:L97").Value = $lArr
$ws.Range("M2:M97").Value = $mArr

Write-Host "Done"